# IWP Review 1 -RECIPE FINDER.pptx
#
# The author's commit fixes the capitalization of the second presenter's
# surname on the title slide: "vaichole" -> "Vaichole", inside the
# subtitle placeholder ("Subtitle 5") text frame, second paragraph
# ("  Tejas vaichole - 19BCE1295").

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

$sh = $null
try {
    $sh = $s.Shapes("Subtitle 5")
} catch {
    $sh = $null
}
if ($sh -eq $null) {
    # Fall back: title slide layout has Rectangle(1), Title(2), Subtitle(3).
    $sh = $s.Shapes.Item(3)
}

$tr = $sh.TextFrame.TextRange

$oldName = "vaichole"
$newName = "Vaichole"

$idx0 = $tr.Text.IndexOf($oldName)
if ($idx0 -ge 0) {
    $chars = $tr.Characters($idx0 + 1, $oldName.Length)
    $chars.Text = $newName
}
